$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 163.63637
$ws.Range("I12").Value = 172.625
$ws.Range("K12").Value = 172.625
$ws.Range("M12").Value = -2.625

$ws.Range("H46").Value = 4672.727
$ws.Range("I46").Value = 3800
$ws.Range("K46").Value = 11400
$ws.Range("M46").Value = -11281

$ws.Range("H60").Value = 4672.727
$ws.Range("I60").Value = 3800
$ws.Range("K60").Value = 11400
$ws.Range("M60").Value = -10916

$ws.Range("H76").Value = 7599
$ws.Range("I76").Value = 7038.6
$ws.Range("K76").Value = 7038.6
$ws.Range("M76").Value = -6723.6

$ws.Range("H79").Value = 7599
$ws.Range("I79").Value = 7038.6
$ws.Range("K79").Value = 7038.6
$ws.Range("M79").Value = -5946.6

$ws.Range("H98").Value = 3963.9697
$ws.Range("J98").Value = 5887.4443
$ws.Range("L98").Value = 5887.4443
$ws.Range("N98").Value = -8883.444299999999

$ws.Range("H107").Value = 2749.1667
$ws.Range("I107").Value = 2500
$ws.Range("K107").Value = 2500
$ws.Range("M107").Value = -580

$ws.Range("H111").Value = 3715.3
$ws.Range("J111").Value = 3159.3333
$ws.Range("L111").Value = 9477.999899999999
$ws.Range("N111").Value = -15611.9999

$ws.Range("H122").Value = 3963.9697
$ws.Range("J122").Value = 5887.4443
$ws.Range("L122").Value = 17662.3329
$ws.Range("N122").Value = -22562.3329

$ws.Range("H125").Value = 3454.8333
$ws.Range("I125").Value = 2267.1428
$ws.Range("J125").Value = 5117.6
$ws.Range("K125").Value = 20404.2852
$ws.Range("L125").Value = 46058.4
$ws.Range("M125").Value = -17944.2852
$ws.Range("N125").Value = -50978.4

$ws.Range("H127").Value = 850.0909
$ws.Range("I127").Value = 785.1
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 2355.3
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 2604.7
$ws.Range("N127").Value = -14420

$ws.Range("H132").Value = 35717740
$ws.Range("J132").Value = 4006
$ws.Range("L132").Value = 12018
$ws.Range("N132").Value = -17078

$ws.Range("H135").Value = 13889718
$ws.Range("I135").Value = 939
$ws.Range("K135").Value = 8451
$ws.Range("M135").Value = -5916

$ws.Range("H138").Value = 5140.541
$ws.Range("I138").Value = 937.24243
$ws.Range("J138").Value = 10094.429
$ws.Range("K138").Value = 2811.72729
$ws.Range("L138").Value = 30283.287
$ws.Range("M138").Value = 2328.27271
$ws.Range("N138").Value = -40563.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 300000
$ws.Range("J34").Value = 500000
$ws.Range("L34").Value = 500000
$ws.Range("N34").Value = -500542

$ws.Range("H60").Value = 71676.664
$ws.Range("I60").Value = 65013
$ws.Range("K60").Value = 65013
$ws.Range("M60").Value = -64280

$ws.Range("H61").Value = 19232516
$ws.Range("I61").Value = 22728982
$ws.Range("K61").Value = 22728982
$ws.Range("M61").Value = -22728770

$ws.Range("H74").Value = 21741692
$ws.Range("J74").Value = 2926.7144
$ws.Range("L74").Value = 2926.7144
$ws.Range("N74").Value = -4674.7144

$ws.Range("H77").Value = 21741692
$ws.Range("J77").Value = 2926.7144
$ws.Range("L77").Value = 14633.572
$ws.Range("N77").Value = -23369.572

$ws.Range("H97").Value = 1660.6
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 20864170
$ws.Range("I132").Value = 1947.7667
$ws.Range("J132").Value = 55634540
$ws.Range("K132").Value = 5843.300099999999
$ws.Range("L132").Value = 166903620
$ws.Range("M132").Value = -3313.300099999999
$ws.Range("N132").Value = -166908680

$ws.Range("H136").Value = 19232516
$ws.Range("I136").Value = 22728982
$ws.Range("K136").Value = 68186946
$ws.Range("M136").Value = -68184396

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3101.25
$ws.Range("I20").Value = 2522.75
$ws.Range("K20").Value = 2522.75
$ws.Range("M20").Value = -2275.75

$ws.Range("H134").Value = 2246.7568
$ws.Range("I134").Value = 2284.4688
$ws.Range("K134").Value = 6853.4064
$ws.Range("M134").Value = -4318.4064

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7199.1333
$ws.Range("I22").Value = 9417.091
$ws.Range("K22").Value = 9417.091
$ws.Range("M22").Value = -9067.091

$ws.Range("H31").Value = 6948079
$ws.Range("I31").Value = 2325.08
$ws.Range("J31").Value = 22733884
$ws.Range("K31").Value = 2325.08
$ws.Range("L31").Value = 22733884
$ws.Range("M31").Value = -2030.08
$ws.Range("N31").Value = -22734474

$ws.Range("H34").Value = 6948079
$ws.Range("I34").Value = 2325.08
$ws.Range("J34").Value = 22733884
$ws.Range("K34").Value = 2325.08
$ws.Range("L34").Value = 22733884
$ws.Range("M34").Value = -2123.08
$ws.Range("N34").Value = -22734288

$ws.Range("H56").Value = 20000
$ws.Range("J56").Value = 20000
$ws.Range("L56").Value = 20000
$ws.Range("N56").Value = -21690

$ws.Range("H58").Value = 3299.5
$ws.Range("I58").Value = 3299.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3299.5
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3096.5

$ws.Range("H68").Value = 95354.875
$ws.Range("J68").Value = 95354.875
$ws.Range("L68").Value = 95354.875
$ws.Range("N68").Value = -96852.875

$ws.Range("H71").Value = 95354.875
$ws.Range("J71").Value = 95354.875
$ws.Range("L71").Value = 286064.625
$ws.Range("N71").Value = -293552.625

$ws.Range("H74").Value = 55000
$ws.Range("J74").Value = 55000
$ws.Range("L74").Value = 55000
$ws.Range("N74").Value = -56748

$ws.Range("H77").Value = 55000
$ws.Range("J77").Value = 55000
$ws.Range("L77").Value = 165000
$ws.Range("N77").Value = -173736

$ws.Range("H107").Value = 694.5238000000001
$ws.Range("I107").Value = 635.6667
$ws.Range("J107").Value = 841.6667
$ws.Range("K107").Value = 635.6667
$ws.Range("L107").Value = 841.6667
$ws.Range("M107").Value = 1284.3333
$ws.Range("N107").Value = -4681.6667

$ws.Range("H132").Value = 2932.8667
$ws.Range("I132").Value = 2855.2222
$ws.Range("K132").Value = 8565.6666
$ws.Range("M132").Value = -6035.6666

$ws.Range("H134").Value = 2461.5
$ws.Range("I134").Value = 2096.75
$ws.Range("K134").Value = 6290.25
$ws.Range("M134").Value = -3755.25

$ws.Range("H136").Value = 3299.5
$ws.Range("I136").Value = 3299.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9898.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -7348.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71.53846
$ws.Range("J2").Value = 221.66667
$ws.Range("L2").Value = 1330.00002
$ws.Range("N2").Value = -1556.00002

$ws.Range("H11").Value = 136.33333
$ws.Range("I11").Value = 104.5
$ws.Range("K11").Value = 313.5
$ws.Range("M11").Value = -173.5

$ws.Range("H129").Value = 1559.9445
$ws.Range("I129").Value = 1025.5714
$ws.Range("K129").Value = 3076.7142
$ws.Range("M129").Value = 1923.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1163.5625
$ws.Range("I97").Value = 1115.4166
$ws.Range("J97").Value = 1308
$ws.Range("K97").Value = 1115.4166
$ws.Range("L97").Value = 1308
$ws.Range("M97").Value = -619.4166
$ws.Range("N97").Value = -2300

$ws.Range("H132").Value = 4994.069
$ws.Range("I132").Value = 4851
$ws.Range("K132").Value = 14553
$ws.Range("M132").Value = -12023

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 542869
$ws.Range("I93").Value = 1753.7826
$ws.Range("K93").Value = 1753.7826
$ws.Range("M93").Value = -505.7826

$ws.Range("H122").Value = 4314579
$ws.Range("I122").Value = 3927.9524
$ws.Range("K122").Value = 11783.8572
$ws.Range("M122").Value = -9333.8572

$ws.Range("H132").Value = 2075.647
$ws.Range("I132").Value = 1469.6072
$ws.Range("K132").Value = 4408.821599999999
$ws.Range("M132").Value = -1878.821599999999

$ws.Range("H136").Value = 1180211.6
$ws.Range("I136").Value = 1336574.1
$ws.Range("K136").Value = 4009722.3
$ws.Range("M136").Value = -4007172.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 46666.668

$ws.Range("H67").Value = 46666.668

$ws.Range("H81").Value = 5723.125
$ws.Range("I81").Value = 3949.6667
$ws.Range("J81").Value = 6787.2
$ws.Range("K81").Value = 7899.3334
$ws.Range("L81").Value = 13574.4
$ws.Range("M81").Value = -6838.3334
$ws.Range("N81").Value = -15696.4

$ws.Range("H84").Value = 5723.125
$ws.Range("I84").Value = 3949.6667
$ws.Range("J84").Value = 6787.2
$ws.Range("K84").Value = 39496.667
$ws.Range("L84").Value = 67872
$ws.Range("M84").Value = -34192.667
$ws.Range("N84").Value = -78480

$ws.Range("H107").Value = 2279.4666
$ws.Range("I107").Value = 1260
$ws.Range("J107").Value = 2789.2
$ws.Range("K107").Value = 3780
$ws.Range("L107").Value = 8367.599999999999
$ws.Range("M107").Value = -1860
$ws.Range("N107").Value = -12207.6

$ws.Range("H132").Value = 1211.0834
$ws.Range("I132").Value = 1139.3636
$ws.Range("K132").Value = 3418.0908
$ws.Range("M132").Value = -888.0907999999999
